# Update the "想去人数" (interested-count) values in the F column
# of both the "展览" sheet and the "全部类型" sheet, reflecting the
# refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" : rows 4,5,7,8,9 in column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 11
$wsExhibit.Range("F5").Value = 3822
$wsExhibit.Range("F7").Value = 50
$wsExhibit.Range("F8").Value = 235
$wsExhibit.Range("F9").Value = 18

# Sheet "全部类型" : rows 8,9,11,13,14 in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 11
$wsAll.Range("F9").Value = 3822
$wsAll.Range("F11").Value = 50
$wsAll.Range("F13").Value = 235
$wsAll.Range("F14").Value = 18
